# Append the latest daily COVID-19 data row to the "Tabela1" table on the
# "Covid-19 podatki" sheet (row for 2020-06-08 / date serial 43990).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row - this keeps the table definition (ref /
# autoFilter) and the sheet's used range in sync, the same way Excel does
# when a user fills in the row right below an Excel Table.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()
$r = $newListRow.Range.Row

# --- Column A (Date) ---------------------------------------------------
# Column A's default format is Text, so the number format has to become a
# real date format *before* the value is written - otherwise the serial
# number would be stored as literal text instead of a date.
$colA = $ws.Cells.Item($r, 1)
$colA.Font.Name = "Calibri Light"
$colA.Font.Size = 10
$colA.HorizontalAlignment = -4152   # xlRight
$colA.VerticalAlignment = -4160     # xlTop
$colA.NumberFormat = "d/ m/ yyyy;@"
$colA.Value2 = 43990

# --- Column B (Tested, all) --------------------------------------------
$colB = $ws.Cells.Item($r, 2)
$colB.Font.Name = "Calibri Light"
$colB.Font.Size = 10
$colB.HorizontalAlignment = -4152   # xlRight
$colB.NumberFormat = "#,##0"
$colB.Value2 = 84130

# --- Columns C-J (remaining daily counters, General number format) -----
$colC = $ws.Cells.Item($r, 3)
$colC.Font.Name = "Calibri Light"
$colC.Font.Size = 10
$colC.HorizontalAlignment = -4152
$colC.Value2 = 814

$colD = $ws.Cells.Item($r, 4)
$colD.Font.Name = "Calibri Light"
$colD.Font.Size = 10
$colD.HorizontalAlignment = -4152
$colD.Value2 = 1486

$colE = $ws.Cells.Item($r, 5)
$colE.Font.Name = "Calibri Light"
$colE.Font.Size = 10
$colE.HorizontalAlignment = -4152
$colE.Value2 = 1

$colF = $ws.Cells.Item($r, 6)
$colF.Font.Name = "Calibri Light"
$colF.Font.Size = 10
$colF.HorizontalAlignment = -4152
$colF.Value2 = 6

$colG = $ws.Cells.Item($r, 7)
$colG.Font.Name = "Calibri Light"
$colG.Font.Size = 10
$colG.HorizontalAlignment = -4152
$colG.Value2 = 0

$colH = $ws.Cells.Item($r, 8)
$colH.Font.Name = "Calibri Light"
$colH.Font.Size = 10
$colH.HorizontalAlignment = -4152
$colH.Value2 = 0

$colI = $ws.Cells.Item($r, 9)
$colI.Font.Name = "Calibri Light"
$colI.Font.Size = 10
$colI.HorizontalAlignment = -4152
$colI.Value2 = 109

$colJ = $ws.Cells.Item($r, 10)
$colJ.Font.Name = "Calibri Light"
$colJ.Font.Size = 10
$colJ.HorizontalAlignment = -4152
$colJ.Value2 = 0

# Keep the cursor where the author's session left it after entering the row.
$ws.Range("E94").Select() | Out-Null
